$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Update commit" / "Update proof" columns (E/F) for the last two
#     benchmark blocks (rows 24-26 and 28-30). ---

# Headers (row 25 / row 29)
$ws.Cells.Item(25,5).Value = "Update commit"
$ws.Cells.Item(25,6).Value = "Update proof"
$ws.Cells.Item(29,5).Value = "Update commit"
$ws.Cells.Item(29,6).Value = "Update proof"

# Data values (row 26 / row 30) - styled with the red font used by the
# other timing values in this workbook.
$ws.Cells.Item(26,5).Value = [double]"1.5449523925781201E-4"
$ws.Cells.Item(26,5).Font.Color = 255
$ws.Cells.Item(26,6).Value = [double]"4.55379486083984E-4"
$ws.Cells.Item(26,6).Font.Color = 255

$ws.Cells.Item(30,5).Value = [double]"4.2510032653808502E-4"
$ws.Cells.Item(30,5).Font.Color = 255
$ws.Cells.Item(30,6).Value = [double]"1.0986328125E-3"
$ws.Cells.Item(30,6).Font.Color = 255

# --- Column widths for the new E/F columns ---
$ws.Range("E1:F1").ColumnWidth = 18.27

# --- Selection / view bookkeeping: the sheet now scrolls back to the top
#     and the active cell moves one row down (D31 -> D32). ---
$ws.Range("D32").Select()
